$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition) updates to column F ("想去人数" / want-to-go count)
$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F2").Value = 4575
$wsExhibition.Range("F3").Value = 2505
$wsExhibition.Range("F4").Value = 481
$wsExhibition.Range("F6").Value = 60
$wsExhibition.Range("F9").Value = 137
$wsExhibition.Range("F10").Value = 177
$wsExhibition.Range("F12").Value = 1712
$wsExhibition.Range("F13").Value = 312
$wsExhibition.Range("F14").Value = 3782
$wsExhibition.Range("F15").Value = 28

# Sheet "全部类型" (All Types) updates to column F ("想去人数" / want-to-go count)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 4575
$wsAll.Range("F3").Value = 2505
$wsAll.Range("F4").Value = 481
$wsAll.Range("F7").Value = 60
$wsAll.Range("F11").Value = 137
$wsAll.Range("F12").Value = 177
$wsAll.Range("F16").Value = 1712
$wsAll.Range("F17").Value = 312
$wsAll.Range("F18").Value = 3782
$wsAll.Range("F19").Value = 28
